$d = $word.ActiveDocument

# ===========================================================================
# Before:
#   ... (paragraphs 1-6: dates/weather, unchanged) ...
#   P7 (last paragraph, has pPr/rPr rFonts hint=eastAsia):
#       run1 "今天天气不错"  (rFonts hint=eastAsia)
#       run2 ",心情也不错"  (rFonts hint=eastAsia)
#
# After:
#   ... (paragraphs 1-6: unchanged) ...
#   new P7: run1 "今天天气不错" + run2 ",心情也不错"   (duplicate of old P7 text)
#   new P8: run1 "星期五"
#   P9 (= old P7, same pPr):
#       run1 "晴，今天学习了分支管理，创建了一个dev分支 使用 Git" (rFonts hint=eastAsia)
#       run2 " "                                                  (no rFonts)
#       run3 "创建分支简单又便捷"                                  (rFonts hint=eastAsia)
# ===========================================================================

# The paragraph we need to duplicate/rewrite is always the last paragraph of
# the body (identified dynamically so this keeps working regardless of the
# exact starting paragraph count).
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$oldText = $lastPara.Range.Text  # includes trailing paragraph mark

# ---------------------------------------------------------------------------
# Step 1: insert a brand-new paragraph right before the last paragraph that
#         duplicates its current text ("今天天气不错" + ",心情也不错"). It is
#         typed as two separate runs, then the two intermediate paragraphs
#         are merged back together (by deleting the paragraph mark between
#         them) so the final paragraph keeps two distinct <w:r> elements,
#         matching how the real edit history produced the target markup.
# ---------------------------------------------------------------------------
$prevPara = $d.Paragraphs.Item($lastIndex - 1)
$prevPara.Range.InsertParagraphAfter()

$dupPara1 = $d.Paragraphs.Item($lastIndex)
$s = $dupPara1.Range.Start
$d.Range($s, $s).Text = "今天天气不错"

$dupPara1 = $d.Paragraphs.Item($lastIndex)
$dupPara1.Range.InsertParagraphAfter()

$dupPara2 = $d.Paragraphs.Item($lastIndex + 1)
$s2 = $dupPara2.Range.Start
$d.Range($s2, $s2).Text = ",心情也不错"

# merge the two paragraphs above back into a single paragraph
$dupEnd = $d.Paragraphs.Item($lastIndex).Range.End
$d.Range($dupEnd - 1, $dupEnd).Delete()

# ---------------------------------------------------------------------------
# Step 2: insert another new paragraph "星期五" right after the duplicate,
#         and before the original (soon to be rewritten) last paragraph.
# ---------------------------------------------------------------------------
$dupPara = $d.Paragraphs.Item($lastIndex)
$dupPara.Range.InsertParagraphAfter()

$weekdayPara = $d.Paragraphs.Item($lastIndex + 1)
$s3 = $weekdayPara.Range.Start
$d.Range($s3, $s3).Text = "星期五"

# ---------------------------------------------------------------------------
# Step 3: rewrite the original last paragraph (now pushed down to index
#         $lastIndex + 2) with the new content:
#   old: "今天天气不错" + ",心情也不错"
#   new: "晴，今天学习了分支管理，创建了一个dev分支 使用 Git" + " " + "创建分支简单又便捷"
# Both Find/Replace calls are scoped to this specific paragraph's range so
# they cannot accidentally match the newly-created duplicate paragraph
# above (which now contains the same original text).
# ---------------------------------------------------------------------------
$targetIndex = $lastIndex + 2
$targetPara = $d.Paragraphs.Item($targetIndex)
$targetRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End)
$targetRange.Find.Execute("今天天气不错", $true, $false, $false, $false, $false, $true, 1, $false, `
    "晴，今天学习了分支管理，创建了一个dev分支 使用 Git", 2) | Out-Null

$targetPara = $d.Paragraphs.Item($targetIndex)
$targetRange = $d.Range($targetPara.Range.Start, $targetPara.Range.End)
$targetRange.Find.Execute(",心情也不错", $true, $false, $false, $false, $false, $true, 1, $false, `
    "创建分支简单又便捷", 2) | Out-Null

# 3c: splice a plain-formatted space (no rFonts at all) in between "Git" and
# "创建分支简单又便捷". Typing a literal space here would normally pick up
# the surrounding rFonts hint="eastAsia" formatting, so instead we copy an
# already-unformatted space character that exists in paragraph 1's text
# ("2022-6-24 ") and paste it at the split point -- paste preserves the
# source formatting (i.e. the lack of it) instead of inheriting from the
# destination.
$p1 = $d.Paragraphs.Item(1)
$spaceSrcStart = $p1.Range.Start + 9
$spaceSrcEnd = $spaceSrcStart + 1
$d.Range($spaceSrcStart, $spaceSrcEnd).Copy()

$targetPara = $d.Paragraphs.Item($targetIndex)
$splitPoint = $targetPara.Range.Start + ([string]"晴，今天学习了分支管理，创建了一个dev分支 使用 Git").Length
$d.Range($splitPoint, $splitPoint).Paste()
